# Add the "Real time (minutes)" (column D) figures that were recorded once
# WalletService.CreateWallet was implemented and its tests were written.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = 30
$ws.Range("D7").Value = 10
$ws.Range("D8").Value = 15
$ws.Range("D9").Value = 45

# Recompute the (wrap-text) row heights now that column D has new content.
$ws.Rows.Item(6).RowHeight = 29.25
$ws.Rows.Item(9).RowHeight = 43.5
$ws.Rows.Item(10).RowHeight = 29.25

# Match the author's final selection in the sheet.
$ws.Range("D10").Select()
